$d = $word.ActiveDocument

$d.Content.Find.Execute("284×9=2556", $true, $false, $false, $false, $false, $true, 1, $false, "327×4=1308", 2) | Out-Null
$d.Content.Find.Execute("771×8=6168", $true, $false, $false, $false, $false, $true, 1, $false, "357×4=1428", 2) | Out-Null
$d.Content.Find.Execute("694×3=2082", $true, $false, $false, $false, $false, $true, 1, $false, "988×4=3952", 2) | Out-Null
$d.Content.Find.Execute("339×6=2034", $true, $false, $false, $false, $false, $true, 1, $false, "142×2=284", 2) | Out-Null
$d.Content.Find.Execute("590×9=5310", $true, $false, $false, $false, $false, $true, 1, $false, "578×8=4624", 2) | Out-Null
$d.Content.Find.Execute("700×6=4200", $true, $false, $false, $false, $false, $true, 1, $false, "753×5=3765", 2) | Out-Null
$d.Content.Find.Execute("794×4=3176", $true, $false, $false, $false, $false, $true, 1, $false, "131×6=786", 2) | Out-Null
$d.Content.Find.Execute("495×9=4455", $true, $false, $false, $false, $false, $true, 1, $false, "214×5=1070", 2) | Out-Null
$d.Content.Find.Execute("453×6=2718", $true, $false, $false, $false, $false, $true, 1, $false, "249×5=1245", 2) | Out-Null
$d.Content.Find.Execute("710×9=6390", $true, $false, $false, $false, $false, $true, 1, $false, "798×4=3192", 2) | Out-Null
$d.Content.Find.Execute("677×4=2708", $true, $false, $false, $false, $false, $true, 1, $false, "470×7=3290", 2) | Out-Null
$d.Content.Find.Execute("219×5=1095", $true, $false, $false, $false, $false, $true, 1, $false, "867×6=5202", 2) | Out-Null
$d.Content.Find.Execute("284×2=568", $true, $false, $false, $false, $false, $true, 1, $false, "954×4=3816", 2) | Out-Null
$d.Content.Find.Execute("846×3=2538", $true, $false, $false, $false, $false, $true, 1, $false, "240×7=1680", 2) | Out-Null
$d.Content.Find.Execute("404×2=808", $true, $false, $false, $false, $false, $true, 1, $false, "389×3=1167", 2) | Out-Null
$d.Content.Find.Execute("638×7=4466", $true, $false, $false, $false, $false, $true, 1, $false, "342×7=2394", 2) | Out-Null
$d.Content.Find.Execute("768×8=6144", $true, $false, $false, $false, $false, $true, 1, $false, "452×2=904", 2) | Out-Null
$d.Content.Find.Execute("583×5=2915", $true, $false, $false, $false, $false, $true, 1, $false, "281×4=1124", 2) | Out-Null
$d.Content.Find.Execute("167×9=1503", $true, $false, $false, $false, $false, $true, 1, $false, "467×7=3269", 2) | Out-Null
$d.Content.Find.Execute("153×2=306", $true, $false, $false, $false, $false, $true, 1, $false, "440×3=1320", 2) | Out-Null
$d.Content.Find.Execute("228×5=1140", $true, $false, $false, $false, $false, $true, 1, $false, "668×4=2672", 2) | Out-Null
$d.Content.Find.Execute("583×8=4664", $true, $false, $false, $false, $false, $true, 1, $false, "512×6=3072", 2) | Out-Null
$d.Content.Find.Execute("674×3=2022", $true, $false, $false, $false, $false, $true, 1, $false, "677×9=6093", 2) | Out-Null
$d.Content.Find.Execute("758×9=6822", $true, $false, $false, $false, $false, $true, 1, $false, "778×8=6224", 2) | Out-Null
$d.Content.Find.Execute("522×3=1566", $true, $false, $false, $false, $false, $true, 1, $false, "806×7=5642", 2) | Out-Null
